# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.741.62"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.634.66"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.639.22"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "1.859.82"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "25.761.40"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.121"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.121.78"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.548"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "1.769.21"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  +2.83%  "

Write-Output "Applied cryptos list refresh"
